$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.453.18'
$ws.Range('E2').Value = '  -3.06%  '

# Row 3
$ws.Range('D3').Value = '1.656.33'
$ws.Range('E3').Value = '  -4.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.00'
$ws.Range('E5').Value = '  -2.34%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.509'
$ws.Range('E6').Value = '  -2.65%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.15'
$ws.Range('E8').Value = '  +0.54%  '

# Row 9
$ws.Range('E9').Value = '  -1.87%  '

# Row 10
$ws.Range('E10').Value = '  -3.01%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0877'
$ws.Range('E11').Value = '  -1.82%  '

# Row 12
$ws.Range('D12').Value = '1.890.67'
$ws.Range('E12').Value = '  -4.26%  '

# Row 13
$ws.Range('D13').Value = '1.654.72'
$ws.Range('E13').Value = '  -4.30%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.12'
$ws.Range('E14').Value = '  -2.88%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.561'
$ws.Range('E15').Value = '  -0.73%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.70'
$ws.Range('E16').Value = '  -2.95%  '

# Row 17
$ws.Range('D17').Value = '27.470.48'
$ws.Range('E17').Value = '  -2.97%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '239.43'
$ws.Range('E18').Value = '  -2.53%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0728'
$ws.Range('E19').Value = '  -3.16%  '

# Row 20
$ws.Range('E20').Value = '  -4.36%  '

# Row 21
$ws.Range('E21').Value = '  +0.00%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.42'
$ws.Range('E22').Value = '  -4.26%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('E23').Value = '  -3.64%  '

# Row 24
$ws.Range('E24').Value = '  -1.48%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.70'
$ws.Range('E25').Value = '  -2.32%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('E26').Value = '  -4.10%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.21'
$ws.Range('E27').Value = '  -2.58%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.11%  '

# Row 29
$ws.Range('E29').Value = '  -2.48%  '

# Row 30
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0498'
$ws.Range('E30').Value = '  -3.84%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.19'
$ws.Range('E31').Value = '  -1.45%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  -3.29%  '

# Row 33
$ws.Range('D33').Value = '1.446.25'
$ws.Range('E33').Value = '  -2.50%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  -5.45%  '

# Row 35
$ws.Range('E35').Value = '  -5.04%  '

# Row 36
$ws.Range('E36').Value = '  -0.90%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.919'
$ws.Range('E37').Value = '  -6.09%  '

# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.572'
$ws.Range('E38').Value = '  -5.37%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0171'
$ws.Range('E39').Value = '  -3.02%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.07'
$ws.Range('E40').Value = '  -0.16%  '

# Row 41
$ws.Range('E41').Value = '  +0.08%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.51'
$ws.Range('E42').Value = '  -4.75%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.43'
$ws.Range('E43').Value = '  -3.75%  '

# Row 44
$ws.Range('E44').Value = '  -3.01%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.789'
$ws.Range('E45').Value = '  -1.63%  '

# Row 46
$ws.Range('D46').Value = '1.799.27'
$ws.Range('E46').Value = '  -4.27%  '

# Row 47
$ws.Range('E47').Value = '  -2.13%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.35'
$ws.Range('E48').Value = '  -2.17%  '

# Row 49
$ws.Range('E49').Value = '  -6.18%  '

# Row 50
$ws.Range('E50').Value = '  -2.11%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.82'
$ws.Range('E51').Value = '  -4.62%  '
